$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '59.420.18'
$ws.Range("E2").Value = '  +0.31%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.527.14'
$ws.Range("E3").Value = '  +2.86%  '

# Row 4
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.96'
$ws.Range("E5").Value = '  +0.79%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.53'
$ws.Range("E6").Value = '  -0.48%  '

# Row 7
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$ws.Range("E8").Value = '  -0.29%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.551.48'
$ws.Range("E9").Value = '  +3.10%  '

# Row 10
$ws.Range("E10").Value = '  +1.00%  '

# Row 11
$ws.Range("E11").Value = '  +0.07%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.60'
$ws.Range("E12").Value = '  +1.47%  '

# Row 13
$ws.Range("E13").Value = '  +3.10%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.973.60'
$ws.Range("E14").Value = '  +2.69%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.86'
$ws.Range("E15").Value = '  -1.65%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '59.317.23'
$ws.Range("E16").Value = '  +0.17%  '

# Row 17
$ws.Range("E17").Value = '  +2.31%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.536.22'
$ws.Range("E18").Value = '  +2.69%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.22'
$ws.Range("E19").Value = '  +0.01%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.30'
$ws.Range("E20").Value = '  -1.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '326.78'
$ws.Range("E21").Value = '  +0.89%  '

# Row 22
$ws.Range("E22").Value = '  +1.70%  '

# Row 23
$ws.Range("E23").Value = '  +2.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.02'
$ws.Range("E24").Value = '  +1.84%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.435'
$ws.Range("E25").Value = '  -3.95%  '

# Row 26
$ws.Range("E26").Value = '  +2.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.992'
$ws.Range("E27").Value = '  +1.15%  '

# Row 28
$ws.Range("E28").Value = '  +2.94%  '

# Row 29
$ws.Range("E29").Value = '  +1.33%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.84'
$ws.Range("E30").Value = '  +1.27%  '

# Row 31
$ws.Range("E31").Value = '  -1.30%  '

# Row 32
$ws.Range("E32").Value = '  -5.89%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.50'
$ws.Range("E33").Value = '  +6.30%  '

# Row 34
$ws.Range("B34").Value = 'Monero'
$ws.Range("C34").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '160.20'
$ws.Range("E34").Value = '  +1.41%  '

# Row 35
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("E35").Value = '  -0.03%  '

# Row 36
$ws.Range("E36").Value = '  +0.25%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.41'
$ws.Range("E37").Value = '  -1.30%  '

# Row 38
$ws.Range("E38").Value = '  -6.15%  '

# Row 39
$ws.Range("E39").Value = '  -3.95%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.99'
$ws.Range("E40").Value = '  +1.46%  '

# Row 41
$ws.Range("B41").Value = 'SuiNetwork'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.834'
$ws.Range("E41").Value = '  +0.22%  '

# Row 42
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.72'
$ws.Range("E42").Value = '  -0.53%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '294.95'
$ws.Range("E43").Value = '  -5.89%  '

# Row 44
$ws.Range("E44").Value = '  -0.25%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.604'
$ws.Range("E45").Value = '  +0.99%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.82'
$ws.Range("E46").Value = '  +0.68%  '

# Row 47
$ws.Range("E47").Value = '  +0.00%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.87'
$ws.Range("E48").Value = '  +1.85%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.40'
$ws.Range("E49").Value = '  -1.81%  '

# Row 50
$ws.Range("E50").Value = '  -1.01%  '

# Row 51
$ws.Range("E51").Value = '  -2.78%  '
